# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") values for the game-log rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$kValues = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    30 = 2
    31 = 0
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
